# Apply "New troops and new abilities" update to allStats workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing row tweaks -------------------------------------------------

# Skeleton (row 13): buff ATK+/Phy DEF+/Magic DEF+
$ws.Range("L13").Value = 10
$ws.Range("M13").Value = 5
$ws.Range("N13").Value = 5

# Demon (row 32): buff ATK+
$ws.Range("L32").Value = 30

# --- New troops ------------------------------------------------------------
# Fill column A (names) first, top to bottom, for every new row so the
# shared-string table gets "Zombie" .. "Skeleton Archer" in row order.
$ws.Range("A34").Value = "Zombie"
$ws.Range("A35").Value = "Skeleton King"
$ws.Range("A36").Value = "Mummy Queen"
$ws.Range("A37").Value = "Mummy"
$ws.Range("A38").Value = "Bat Witch"
$ws.Range("A39").Value = "Skeleton Archer"

# Multi-ability strings are entered in this order: row35, row39, then row36.
$ws.Range("J35").Value = "29, 9, 11"
$ws.Range("J39").Value = "9, 11, 12"
$ws.Range("J36").Value = "11, 30"

# Row 34: Zombie
$ws.Range("B34").Value = "-"
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 40
$ws.Range("E34").Value = 30
$ws.Range("F34").Value = "Phy"
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 2
$ws.Range("J34").Value = 11
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 1
$ws.Range("Q34").Value = 1

# Row 35: Skeleton King
$ws.Range("B35").Value = 1600
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 100
$ws.Range("E35").Value = 60
$ws.Range("F35").Value = "Phy"
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 30
$ws.Range("I35").Value = 2
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 5
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 5
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 1
$ws.Range("Q35").Value = 2

# Row 36: Mummy Queen
$ws.Range("B36").Value = 800
$ws.Range("C36").Value = 4
$ws.Range("D36").Value = 300
$ws.Range("E36").Value = 10
$ws.Range("F36").Value = "Magic"
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 2
$ws.Range("K36").Value = 20
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 1
$ws.Range("Q36").Value = 1

# Row 37: Mummy
$ws.Range("B37").Value = "-"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 150
$ws.Range("E37").Value = 30
$ws.Range("F37").Value = "Phy"
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 3
$ws.Range("J37").Value = 11
$ws.Range("K37").Value = 10
$ws.Range("L37").Value = 10
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 1
$ws.Range("Q37").Value = 1

# Row 38: Bat Witch
$ws.Range("B38").Value = 1200
$ws.Range("C38").Value = 3
$ws.Range("D38").Value = 100
$ws.Range("E38").Value = 30
$ws.Range("F38").Value = "Magic"
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 10
$ws.Range("I38").Value = 3
$ws.Range("J38").Value = 31
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 5
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 10
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 3
$ws.Range("Q38").Value = 4

# Row 39: Skeleton Archer
$ws.Range("B39").Value = 300
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 50
$ws.Range("E39").Value = 60
$ws.Range("F39").Value = "Phy"
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 2
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 3
$ws.Range("Q39").Value = 3

# Row 40: mark the end of the new block with an underlined (empty) cell.
$ws.Range("B40").Font.Underline = $true

# Leave the selection where the author left it when they saved.
$ws.Range("L32").Select()

Write-Output "edit applied"
